# Build site at 2023-04-12 14:53:07 UTC
# Fix Objetivos/Programa resumido/Programa/Bibliografia content, add two
# "Docentes responsaveis" rows (Danubia + Robson), and add the new PT
# objectives / short syllabus / syllabus / bibliography text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the "Objetivos:" row (row 10) - it currently (wrongly) holds the
#    first docente's name; replace with the actual PT objectives text.
$ws.Range("B10:C10").Value = @'
O curso tem como objetivos: proporcionar aos discentes conhecimentos básicos sobre as formas do relevo e o papel da água como agente geomorfológico; Estudar o sistema básico de circulação de água e os padrões de escoamento; Estudar o sistema fluvial sob a perspectiva da análise ambiental do meio físico enfatizando processos geomorfológicos e hidrológicos; Fornecer ao aluno os conceitos básicos de hidrologia aplicadas ao meio ambiente.
'@

# 2) Insert two new rows right after "Docentes responsaveis:" (row 12) to
#    hold the two docente names that were previously misplaced further down.
$ws.Rows("13:14").Insert()

$ws.Range("B13").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C13").Value = "9146830 - Danúbia Caporusso Bargos"

$ws.Range("B14").Value = "7455355 - Robson da Silva Rocha"
$ws.Range("C14").Value = "7455355 - Robson da Silva Rocha"

# 3) "Programa resumido:" row (now row 15) - replace the wrong date value
#    with the actual PT short-syllabus text.
$ws.Range("B15:C15").Value = @'
Geomorfologia Fluvial; Padrões de Drenagem; Escoamentos hidráulicos; medidores; bocais; instrumentos de medição
'@

# 4) "Programa:" row (now row 17) - replace the wrong value (it had picked
#    up the first docente's name) with the actual PT full syllabus text.
$ws.Range("B17:C17").Value = @'
- As teorias geomorfológicas;- Processos e Formas do relevo;- Processos fluviais, morfologias fluviais e padrões de drenagem;- Precipitação;- Infiltração;- Evapotranspiração;- Escoamento superficial;- Instrumentos de medição (Calhas, vertedores e registros);- Operação de reservatórios;- Vazões máximas e mínimas: distribuição de frequência, hidrograma unitário.- Água subterrânea, aquíferos e poços;
'@

# 5) "Bibliografia:" row (now row 23) - replace the recovery-exam text
#    (which had been shifted up one row) with the real bibliography.
$ws.Range("B23:C23").Value = @'
Barth, F.T. et al. - Modelos para Gerenciamento de Recursos Hídricos. São Paulo: Nobel: ABRH (Coleção ABRH de Recursos Hídricos, Vol. 1)., 1987.Pinto, N.L.S. et al. - Hidrologia Básica. São Paulo: Edgard Blucher, 1976.Tucci, C.E.M. - Hidrologia: Ciência e Aplicação, Porto Alegre, Editora da Universidade - ABRH - EPUSP, 1993, (coleção ABRH de Recursos Hídricos, Vol. 4).Villela, S.M. e Mattos. - Hidrologia Aplicada, São Paulo: Mc Graw-hill do Brasil, 1975.Wilson - Engineering Hydrology, London: Mcmillan, 1969.Roberto, A. N., Porto. R.L.L. e Zahed, K.F. - Sistema de Suporte a Decisões para Análise de Cheias em Bacias Complexas. Anais do XII Simpósio Brasileiro de Recursos Hídricos, 1997.Tucci, C.E.M., Porto, R.L.L. e Barros, M.T. - Drenagem Urbana, Porto Alegre, Editora da Universidade - ABRH - UFRGS, 1995, (coleção ABRH de Recursos Hídricos, Vol. 5)Wanielista, M.,Kersten, R. e Eaglin,R. -Hydrology - Water Quantity and Quality Control, John Wiley & Sons, Inc., 567 pág., 1997.Porto, R.L.L. - Técnicas Quantitativas para o Gerenciamento de Recursos Hídricos, Porto Alegre, Editora da Universidade - ABRH - UFRGS, 1997, (coleção ABRH de Recursos Hídricos, Vol. 6)Ward, A.D. Trimble, S.W. – Environmental Hydrology, 2004, Lewis Publishers, 462 ppBrutsaert, W. – Hydrology: An Introduction, Cambridge University Press, 618 p, ISBN 0521824796, 2005.Gordon, N.D McMahon, T.A. Finlayson, B.L. Gippel, C.J. Nathan, R.J. – Stream Hydrology: An Introduction to Ecologists, Second Edition, John Wiley & Sons Ltd., UK, 526 pp.Brooks,K.N. Ffolliott,P.F. Gregersen,H.M. DeBano,L.F. – Hydrology and the Management of Watersheds, Iowa State University Press, 574 pp, 2003Ghosh,S.N. Desai,V.R. – Environmental Hydrology And Hydraulics: Eco-technological Practices for Sustainable Development, Science Publishers, 416 p, ISBN 978-1-57808-403-6, 2006Shaw, E.M. – Hydrology in Practice, Van Nostrand Rainhold, 2007Baird, A.J. Wilby, R.L. – Eco-Hydrology: Plants and water in terrestrial and aquatic environments, Routledge, 2007CHRISTOFOLETTI, A. Geomorfologia Fluvial. São Paulo: Edgar Blucher Ltda, 1981. 313 p. CHIOSSI, N. Geologia de Engenharia. São Paulo: Oficina de Textos, 2013.WICANDER, R.; MONROE, J.S. Geologia. São Paulo: Cengage Learning, 2017.SILVA, L.P. Hidrologia: Engenharia e meio ambiente. Rio de Janeiro: Elsevier, 2015.POLETO, C. Bacias hidrográficas e recursos hídricos. Rio de Janeiro: Interciência, 2014.
'@

# Row-height bookkeeping: rows 13/14 are plain rows (no A label, default
# height); row 17 grows to match the longer syllabus text while row 15
# (the short blurb) stays compact.
$ws.Rows("13:14").RowHeight = 15
$ws.Rows(15).RowHeight = 60
$ws.Rows(17).RowHeight = 120

# Column layout cleanup: column A only needs the label width; column B
# keeps the wider content width (column C already has its own definition).
$ws.Columns("A").ColumnWidth = 30.7109375
$ws.Columns("B").ColumnWidth = 60.7109375
